# "Structure of two diameter" — split the existing Conduit 1 row into two
# stacked pull rows (row 2 + new row 3), merging the cells that remain
# shared across both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (first pull) with its new figures -------------------
$ws.Range("D2").Value = 2
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 3.5

# H2 / J2 hold percentage-looking text (not real numbers) in the source
# file, so force text entry with a leading apostrophe, then restore the
# plain (non quote-prefixed) number format/style by pasting formats from
# a sibling cell that already carries the desired style.
$ws.Range("H2").Value = "'28.07%"
$ws.Range("J2").Value = "'20.62%"

$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Add new row 3 (second pull for the same conduit) ------------------
# Clone row 2's formatting down to row 3 first, then fill in only the
# cells that differ for the second pull (D3, E3); the rest stay blank.
$ws.Range("A2:J2").Copy() | Out-Null
$ws.Range("A3:J3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "7C#14"

# --- Merge the columns shared between the two pulls of Conduit 1 -------
$ws.Range("A2:A3").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("C2:C3").Merge()
$ws.Range("F2:F3").Merge()
$ws.Range("G2:G3").Merge()
$ws.Range("H2:H3").Merge()
$ws.Range("I2:I3").Merge()
$ws.Range("J2:J3").Merge()
